# Updated symbol list on Tue Dec 27 07:52:31 UTC 2022 with GitHub Actions
#
# Applies the "Price" (column D) refresh, a couple of label/volume-tag
# corrections (column E), and the CEJI / BKEXToken row swap (columns B-E,
# rows 42-43) coming from the upstream coinranking.com scrape.
#
# All of these cells store their content as text even though many of them
# look like plain numbers (e.g. "243.26"), so a plain `.Value = ...` would
# let Excel's COM layer re-interpret the string as a numeric literal and
# silently change the cell's type/style. We briefly force NumberFormat to
# Text ("@") so the assignment is kept as a string, then restore the
# "Normal" style so the cell's formatting ends up exactly as it started
# (only the stored text content changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Value
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
    $rng.Style = "Normal"
}

# --- Column D (Price) refresh -------------------------------------------
Set-TextValue "D2"  "243.26"
Set-TextValue "D3"  "23.03"
Set-TextValue "D4"  "5.409"
Set-TextValue "D5"  "0.05981"
Set-TextValue "D6"  "3.423"
Set-TextValue "D7"  "6.499"
Set-TextValue "D8"  "0.8126"
Set-TextValue "D9"  "0.9308"
Set-TextValue "D10" "0.1443"
Set-TextValue "D11" "0.07453"
Set-TextValue "D13" "0.03054"
Set-TextValue "D14" "0.09351"
Set-TextValue "D15" "3.855"
Set-TextValue "D16" "0.001580"
Set-TextValue "D17" "0.04703"
Set-TextValue "D18" "0.0005943"
Set-TextValue "D19" "0.005902"
Set-TextValue "D20" "0.001269"
Set-TextValue "D21" "0.004876"
Set-TextValue "D24" "2.133"
Set-TextValue "D25" "0.3238"
Set-TextValue "D26" "0.1332"
Set-TextValue "D27" "0.0002341"
Set-TextValue "D40" "0.03959"
Set-TextValue "D41" "0.006331"
Set-TextValue "D44" "0.008916"
Set-TextValue "D45" "0.00005175"
Set-TextValue "D47" "0.6703"

# --- Column E label corrections ------------------------------------------
Set-TextValue "E20" "19BitKanKANBestin24h"

# --- Rows 42/43: CEJI and BKEXToken swapped rank ---------------------------
Set-TextValue "B42" "BKEXToken"
Set-TextValue "C42" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1076"
Set-TextValue "E42" "41BKEXTokenBKK"

Set-TextValue "B43" "CEJI"
Set-TextValue "C43" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002651"
Set-TextValue "E43" "42CEJICEJI"
